$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = -12.0269
$ws.Range("C13").Value = -12.29499999999999
$ws.Range("C16").Value = -11.9768
$ws.Range("C18").Value = -14.09679999999998
$ws.Range("C20").Value = -13.34489999999998
